# The commit removes a single row from the "posts" sheet (the post about
# "「サイのざらざらした皮を利用して自分の体を掻くネコ」" that used to live at row 769).
# Deleting that row shifts every subsequent row up by one, shrinking the
# used range from A1:C836 down to A1:C835 - exactly what the diff shows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(769).Delete()
